$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that simply had price/volume updated ---
$ws.Range('D2').Value = '51.674.26'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.082.09'
$ws.Range('E3').Value = '  +3.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '388.00'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.546'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.05'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').Value = '3.567.47'
$ws.Range('E13').Value = '  +3.30%  '
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.83'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '3.081.49'
$ws.Range('E16').Value = '  +2.94%  '
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.74'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').Value = '51.769.48'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('E20').Value = '  +2.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.50'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.35'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.02'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('E25').Value = '  -2.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.23'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.96'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.172'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '35.18'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('E34').Value = '  -2.40%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  +1.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.297'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +9.61%  '
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.96'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.25%  '
$ws.Range('D49').Value = '2.041.19'
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('D50').Value = '3.382.55'
$ws.Range('E50').Value = '  +2.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.209'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.21%  '

# --- Rows 35/36 and 47/48 swapped rank order (coin + link + price + volume) ---
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0452'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.81%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '50.18'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.65%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.09'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.10%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.46'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.92%  '
